$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Due Dt" column (G) ---

# Header
$ws.Range("G1").Value = "Due Dt"

# Data rows: G2:G16 (Due Dt values, as Excel date serials)
$dueDates = @{
    2  = 45536
    3  = 45536
    4  = 45536
    5  = 45536
    6  = 45536
    7  = 45505
    8  = 45505
    9  = 45536
    10 = 45536
    11 = 45505
    12 = 45505
    13 = 45536
    14 = 45536
    15 = 45536
    16 = 45536
}

foreach ($row in 2..16) {
    $ws.Range("G$row").Value = $dueDates[$row]
    # Copy the existing date-formatted style from column F so the new
    # column matches the "Completed Dt" number formatting (style index 1).
    $ws.Range("F$row").Copy() | Out-Null
    $ws.Range("G$row").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false

# Row 6's "Completed Dt" (F6) value is removed, leaving just the style.
$ws.Range("F6").ClearContents()

# Update the selection to reflect the last edited cell.
$ws.Range("G16").Select() | Out-Null

# Restore a plain page setup (portrait orientation) for the sheet.
$ws.PageSetup.Orientation = 1
